# daily auto push: 2026-01-16 06:50 UTC
# Insert a new data row for 2026/01/16 14:00 (金) just before the old row 661
# (2026/12/29 火 13:00), pushing all the following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 661..702 down to 662..703, leaving row 661 blank for the new entry.
$ws.Rows.Item(661).Insert()

$dateCell = $ws.Cells.Item(661, 1)
# Force the date column to be stored as literal text (matching every other
# row in the sheet) instead of letting Excel auto-convert the "yyyy/mm/dd"
# looking string into a real date serial number.
$dateCell.NumberFormat = "@"
$dateCell.Value2 = "2026/01/16"
$dateCell.ClearFormats()

$ws.Cells.Item(661, 2).Value = "金"
$ws.Cells.Item(661, 3).Value = 14
$ws.Cells.Item(661, 4).Value = 201
